$d = $word.ActiveDocument

# The "AbstractFirstParagraph" paragraph style inherits a 720-twip
# first-line indent from its base style "Abstract". Fix the abstract
# indent by overriding FirstLineIndent to 0 directly on
# "AbstractFirstParagraph" (adds a <w:pPr><w:ind w:firstLine="0"/></w:pPr>
# to the style definition in styles.xml).
$style = $d.Styles("AbstractFirstParagraph")
$style.ParagraphFormat.FirstLineIndent = 0
